# Add a new "Sheet2" after the existing "Sheet1" — this is the new
# "credit cost / stage-attempts" calculator described in the commit
# message.
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Raw per-tier credit costs, entered in columns D (credits) and E
# (credits for the harder variant).
$dData = @(3200, 6400, 12800, 32000, 64000, 128000)
$eData = @(4000, 8000, 16000, 40000, 80000, 160000)

for ($i = 0; $i -lt 6; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 4).Value = $dData[$i]
    $ws2.Cells.Item($row, 5).Value = $eData[$i]
}

# Match the existing "centered" number style used throughout Sheet1
# (style index 1) instead of inventing a brand-new named style.
$ws1.Range("A1").Copy()
$ws2.Range("D1:E6").PasteSpecial(-4122)  # xlPasteFormats

# Columns A/B accumulate the running total cost, read off of D/E.
$ws2.Range("A1").Formula = "=0+D1"
$ws2.Range("B1").Formula = "=0+E1"
$ws2.Range("A2").Formula = "=A1+D2"
$ws2.Range("B2").Formula = "=B1+E2"
$ws2.Range("A3:B6").Formula = "=A2+D3"

# Sheet2 becomes the active tab/selection, matching the saved view state.
$ws2.Activate() | Out-Null
$ws2.Range("F7").Select() | Out-Null
